$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 329.33334
$ws.Range("I2").Value = 252.85715
$ws.Range("J2").Value = 396.25
$ws.Range("K2").Value = 252.85715
$ws.Range("L2").Value = 396.25
$ws.Range("M2").Value = -139.85715
$ws.Range("N2").Value = -622.25
$ws.Range("H15").Value = 589.38
$ws.Range("I15").Value = 589.38
$ws.Range("K15").Value = 1768.14
$ws.Range("M15").Value = -1599.14
$ws.Range("H17").Value = 3874.5789
$ws.Range("J17").Value = 3323
$ws.Range("L17").Value = 9969
$ws.Range("N17").Value = -10305
$ws.Range("H64").Value = 3474.8333
$ws.Range("J64").Value = 3569.8
$ws.Range("L64").Value = 3569.8
$ws.Range("N64").Value = -4065.8
$ws.Range("H67").Value = 3474.8333
$ws.Range("J67").Value = 3569.8
$ws.Range("L67").Value = 3569.8
$ws.Range("N67").Value = -5285.8
$ws.Range("H116").Value = 21850
$ws.Range("I116").Value = 51155.5
$ws.Range("J116").Value = 2313
$ws.Range("K116").Value = 51155.5
$ws.Range("L116").Value = 2313
$ws.Range("M116").Value = -47713.5
$ws.Range("N116").Value = -9197
$ws.Range("H138").Value = 3155.5833
$ws.Range("I138").Value = 3434.4075
$ws.Range("J138").Value = 3046.4783
$ws.Range("K138").Value = 10303.2225
$ws.Range("L138").Value = 9139.4349
$ws.Range("M138").Value = -5163.2225
$ws.Range("N138").Value = -19419.4349

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16196.403
$ws.Range("I32").Value = 12333.462
$ws.Range("K32").Value = 12333.462
$ws.Range("M32").Value = -12046.462
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").Value = ""
$ws.Range("H74").Value = 969.5454999999999
$ws.Range("I74").Value = 766.575
$ws.Range("J74").Value = 2999.25
$ws.Range("K74").Value = 766.575
$ws.Range("L74").Value = 2999.25
$ws.Range("M74").Value = 107.425
$ws.Range("N74").Value = -4747.25
$ws.Range("H77").Value = 969.5454999999999
$ws.Range("I77").Value = 766.575
$ws.Range("J77").Value = 2999.25
$ws.Range("K77").Value = 3832.875
$ws.Range("L77").Value = 14996.25
$ws.Range("M77").Value = 535.125
$ws.Range("N77").Value = -23732.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1001000
$ws.Range("H89").Value = 1001000
$ws.Range("H99").Value = 933.3333
$ws.Range("I99").Value = 933.3333
$ws.Range("K99").Value = 933.3333
$ws.Range("M99").Value = 564.6667
$ws.Range("H105").Value = 2705.45
$ws.Range("I105").Value = 2479.4211
$ws.Range("K105").Value = 2479.4211
$ws.Range("M105").Value = -732.4211
$ws.Range("H125").Value = 22000
$ws.Range("J125").Value = 22000
$ws.Range("L125").Value = 22000
$ws.Range("N125").Value = -31840

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1918.4
$ws.Range("I16").Value = 1918.4
$ws.Range("K16").Value = 1918.4
$ws.Range("M16").Value = -1631.4
$ws.Range("H31").Value = 5019.4736
$ws.Range("I31").Value = 2200
$ws.Range("J31").Value = 5351.1763
$ws.Range("K31").Value = 2200
$ws.Range("L31").Value = 5351.1763
$ws.Range("M31").Value = -1905
$ws.Range("N31").Value = -5941.1763
$ws.Range("H34").Value = 5019.4736
$ws.Range("I34").Value = 2200
$ws.Range("J34").Value = 5351.1763
$ws.Range("K34").Value = 2200
$ws.Range("L34").Value = 5351.1763
$ws.Range("M34").Value = -1998
$ws.Range("N34").Value = -5755.1763
$ws.Range("H58").Value = 967954.25
$ws.Range("I58").Value = 4349846
$ws.Range("K58").Value = 4349846
$ws.Range("M58").Value = -4349643
$ws.Range("H62").Value = 3766.6667
$ws.Range("I62").Value = 3366.6667
$ws.Range("K62").Value = 3366.6667
$ws.Range("M62").Value = -2742.6667
$ws.Range("H65").Value = 3766.6667
$ws.Range("I65").Value = 3366.6667
$ws.Range("K65").Value = 16833.3335
$ws.Range("M65").Value = -13713.3335
$ws.Range("H86").Value = 58831588
$ws.Range("I86").Value = 111113090
$ws.Range("K86").Value = 111113090
$ws.Range("M86").Value = -111111967
$ws.Range("H89").Value = 58831588
$ws.Range("I89").Value = 111113090
$ws.Range("K89").Value = 555565450
$ws.Range("M89").Value = -555559834
$ws.Range("H105").Value = 986.5714
$ws.Range("I105").Value = 1086.25
$ws.Range("K105").Value = 1086.25
$ws.Range("M105").Value = 660.75
$ws.Range("H107").Value = 2499.7273
$ws.Range("I107").Value = 2338.611
$ws.Range("K107").Value = 2338.611
$ws.Range("M107").Value = -418.6109999999999
$ws.Range("H113").Value = 1918.4
$ws.Range("I113").Value = 1918.4
$ws.Range("K113").Value = 1918.4
$ws.Range("M113").Value = 251.5999999999999
$ws.Range("H135").Value = 52874.5
$ws.Range("J135").Value = 52874.5
$ws.Range("L135").Value = 52874.5
$ws.Range("N135").Value = -63014.5
$ws.Range("H136").Value = 967954.25
$ws.Range("I136").Value = 4349846
$ws.Range("K136").Value = 13049538
$ws.Range("M136").Value = -13046988
$ws.Range("H141").Value = 79819.3
$ws.Range("J141").Value = 80688.11
$ws.Range("L141").Value = 80688.11
$ws.Range("N141").Value = -91048.11

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 18904.092
$ws.Range("I131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("M131").Value = ""
$ws.Range("H134").Value = 34462.97
$ws.Range("I134").Value = 41935.4
$ws.Range("J134").Value = 3327.8333
$ws.Range("K134").Value = 125806.2
$ws.Range("L134").Value = 9983.499899999999
$ws.Range("M134").Value = -120736.2
$ws.Range("N134").Value = -20123.4999
$ws.Range("H137").Value = 4046.8696
$ws.Range("I137").Value = 1904.7693
$ws.Range("K137").Value = 5714.3079
$ws.Range("M137").Value = -614.3078999999998
$ws.Range("H138").Value = 2099.8462
$ws.Range("I138").Value = 2099.8462
$ws.Range("K138").Value = 6299.5386
$ws.Range("M138").Value = -1159.5386
$ws.Range("H140").Value = 2312.7576
$ws.Range("J140").Value = 3543.6667
$ws.Range("L140").Value = 10631.0001
$ws.Range("N140").Value = -20991.0001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5858
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 5858
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 5858
$ws.Range("M70").Value = ""
$ws.Range("N70").Value = -6398
$ws.Range("H73").Value = 5858
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 5858
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 5858
$ws.Range("M73").Value = ""
$ws.Range("N73").Value = -7730
$ws.Range("H80").Value = 3665.6667
$ws.Range("I80").Value = 3665.6667
$ws.Range("K80").Value = 3665.6667
$ws.Range("M80").Value = -2667.6667
$ws.Range("H83").Value = 3665.6667
$ws.Range("I83").Value = 3665.6667
$ws.Range("K83").Value = 18328.3335
$ws.Range("M83").Value = -13336.3335
$ws.Range("H126").Value = 2461959.8
$ws.Range("I126").Value = 3475503.5
$ws.Range("K126").Value = 10426510.5
$ws.Range("M126").Value = -10424040.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2218.4666
$ws.Range("J46").Value = 2509.75
$ws.Range("L46").Value = 2509.75
$ws.Range("N46").Value = -2885.75
$ws.Range("H47").Value = 1000000000
$ws.Range("I47").Value = 1000000000
$ws.Range("K47").Value = 1000000000
$ws.Range("M47").Value = -999999510
$ws.Range("H52").Value = 1000000000
$ws.Range("I52").Value = 1000000000
$ws.Range("K52").Value = 1000000000
$ws.Range("M52").Value = -999999767

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 19980
$ws.Range("J40").Value = 19980
$ws.Range("L40").Value = 19980
$ws.Range("N40").Value = -20278
$ws.Range("H107").Value = 773.52
$ws.Range("I107").Value = 606.7273
$ws.Range("J107").Value = 1996.6666
$ws.Range("K107").Value = 1820.1819
$ws.Range("L107").Value = 5989.9998
$ws.Range("M107").Value = 99.81809999999996
$ws.Range("N107").Value = -9829.9998
$ws.Range("H113").Value = 1244.75
$ws.Range("I113").Value = 989.5
$ws.Range("K113").Value = 2968.5
$ws.Range("M113").Value = -798.5
$ws.Range("H126").Value = 2404.3462
$ws.Range("I126").Value = 2335.5652
$ws.Range("J126").Value = 2931.6667
$ws.Range("K126").Value = 7006.6956
$ws.Range("L126").Value = 8795.000100000001
$ws.Range("M126").Value = -4536.6956
$ws.Range("N126").Value = -13735.0001
$ws.Range("H136").Value = 25254762
$ws.Range("I136").Value = 42736544
$ws.Range("J136").Value = 3300.5557
$ws.Range("K136").Value = 128209632
$ws.Range("L136").Value = 9901.667099999999
$ws.Range("M136").Value = -128207082
$ws.Range("N136").Value = -15001.6671
